# WebForm User Assignment execution
#
# Fills in the PN_Value (phone number) column F for rows 2-18 with the
# values assigned to each user during this execution run, and flips the
# AN2 "Match1UserPos" flag to match AO2's "Match2UserPos" value (both
# become "2") now that a second match was found.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($a1, $value) {
    # Leading apostrophe forces Excel to store the numeric-looking string
    # as literal text (same as typing it in manually), then resetting the
    # cell style back to Normal clears the quote-prefix formatting that
    # entry would otherwise leave behind.
    $cell = $ws.Range($a1)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue "F2"  "9840072974"
Set-TextValue "F3"  "9840067309"
Set-TextValue "F4"  "9840004187"
Set-TextValue "F5"  "9840015195"
Set-TextValue "F6"  "9840064524"
Set-TextValue "F7"  "9840057900"
Set-TextValue "F8"  "9840063447"
Set-TextValue "F9"  "9840043010"
Set-TextValue "F10" "9840067881"
Set-TextValue "F11" "9840012679"
Set-TextValue "F12" "9840089552"
Set-TextValue "F13" "9840077695"
Set-TextValue "F14" "9840001205"
Set-TextValue "F15" "9840003702"
Set-TextValue "F16" "9840058190"
Set-TextValue "F17" "9840028942"
Set-TextValue "F18" "9840052972"

Set-TextValue "AN2" "2"
